# EWD-23058 - External review: As a paid feature
# Adds two new localization keys ("reviewNotPaidMessage" and
# "upgradeToStarterPlanToUseCommentsErrorMessage") right after the existing
# "reviewNoComments" / "No comments yet" row on the main Sheet1, pushing the
# "Feedback" block (and everything after it) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two fresh rows right before the "Feedback" row (row 360) - this
# shifts every subsequent row (incl. the trailing formatting-only row) down
# by two, matching the target layout.
$ws.Rows("359:360").Insert()

# New row 359: reviewNotPaidMessage
$ws.Range("B359").Value = "reviewNotPaidMessage"
$ws.Range("C359").Value = 'Starter plan users only. Please <u style="color: #2cb7ff;">upgrade</u> to use this functionality.'
$ws.Range("D359").Value = 'Starter plan users only. Please <u style="color: #2cb7ff;">upgrade</u> to use this functionality.'
$ws.Range("E359").Value = 'Starter plan users only. Please <u style="color: #2cb7ff;">upgrade</u> to use this functionality.'

# New row 360: upgradeToStarterPlanToUseCommentsErrorMessage
$ws.Range("B360").Value = "upgradeToStarterPlanToUseCommentsErrorMessage"
$ws.Range("C360").Value = 'Please upgrade your account to "Starter plan" to be able to see comments.'
$ws.Range("D360").Value = 'Please upgrade your account to "Starter plan" to be able to see comments.'
$ws.Range("E360").Value = 'Please upgrade your account to "Starter plan" to be able to see comments.'

# C360 picks up the "translated" (green) status formatting, matching the
# style used elsewhere in the sheet (e.g. the EN column cells of the
# surrounding rows).
$ws.Range("C363").Copy()
$ws.Range("C360").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "Applied EWD-23058 edit: inserted reviewNotPaidMessage / upgradeToStarterPlanToUseCommentsErrorMessage rows."
